# Sweden Allsvenskan workbook update
# The underlying data source was re-fetched and a handful of match rows
# ended up re-ordered relative to each other (their "id" in column B, and
# all associated odds data in columns B:AC, got shuffled between rows that
# share the same match date). This script restores each row's correct
# B:AC payload by redistributing the values among the affected rows,
# leaving column A (the sequential row index) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($Row1, $Row2) {
    $range1 = $ws.Range("B$Row1`:AC$Row1")
    $range2 = $ws.Range("B$Row2`:AC$Row2")

    $data1 = $range1.Value2
    $data2 = $range2.Value2

    $range1.Value2 = $data2
    $range2.Value2 = $data1
}

function Cycle-RowData($Rows) {
    $count = $Rows.Length

    # Snapshot all of the original B:AC payloads first.
    $originals = @()
    foreach ($r in $Rows) {
        $originals += , ($ws.Range("B$r`:AC$r").Value2)
    }

    # Each row takes on the payload that originally belonged to the next
    # row in the list (wrapping around at the end).
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $Rows[$i]
        $sourceData = $originals[($i + 1) % $count]
        $ws.Range("B$targetRow`:AC$targetRow").Value2 = $sourceData
    }
}

# Simple pairwise swaps (two rows that share a date exchanged their data).
Swap-RowData 47  48
Swap-RowData 52  53
Swap-RowData 56  57
Swap-RowData 67  68
Swap-RowData 88  89
Swap-RowData 103 104
Swap-RowData 165 166
Swap-RowData 172 174
Swap-RowData 196 197
Swap-RowData 228 229

# A 5-row group that rotates rather than simply swapping pairwise.
Cycle-RowData @(235, 237, 239, 240, 238)
